$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.758.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.289.35'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.80%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '96.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.25'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0936'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.92'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("E13").Value = '  +1.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.75'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.632.71'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.297.79'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.719.09'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("E19").Value = '  +4.15%  '
$ws.Range("E20").Value = '  -2.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +12.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.60'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.49%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '38.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.12%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '174.99'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0894'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  -0.36%  '
$ws.Range("E35").Value = '  +0.48%  '
$ws.Range("E36").Value = '  +6.34%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  -1.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.79%  '
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.32'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.25'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.12%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.19'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.90%  '
$ws.Range("B44").Value = 'ARBITRUM'
$ws.Range("C44").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.39%  '
$ws.Range("E45").Value = '  -3.04%  '
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("E47").Value = '  -0.30%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.97'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.50%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.440'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.25%  '
